# Update the "top-level-description" cells in column B with revised, shorter
# descriptive text, as part of the "Latest generated outputs 2025-09-10"
# regeneration of the permission-in-principle-pip specification sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Specification")

$ws.Range("B21").Value = "Name and contact information if an agent is being used."
$ws.Range("B29").Value = "Name and contact information if an agent is being used."
$ws.Range("B33").Value = "Name and contact information for the parties making the application."
$ws.Range("B39").Value = "Telephone number and email address of the applicant."
$ws.Range("B43").Value = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
$ws.Range("B44").Value = "Details of any conflict of interest that may exist between the applicant and planning authority."
$ws.Range("B47").Value = "Signed and dated verification of the application's accuracy."
$ws.Range("B50").Value = "Details of the residential and non-residential parts of the proposed development."
$ws.Range("B57").Value = "Where the proposed development will be built."
$ws.Range("B66").Value = "Any additional relevant information about the development site."

$wb.Save()
